# UserAttributes.xlsx - "SQL edits and Dashboard Personal updates"
#
# Sheet1 is a matrix of attribute-name rows (A..F = UserModel, UserModel_Slim,
# Dashboard, SearchFor, ViewAllAtt, UserInformation) marking which DTO/view
# exposes which property. This change:
#   1. Adds "AgeRange" to the UserInformation column (F3) - row 3 is the
#      AgeRange row.
#   2. Moves "ProfilePicPath" off the Dashboard column in row 20 (C20 cleared)
#      and re-adds it one row down in a brand-new row 21, this time exposing
#      it on both UserModel (A21) and Dashboard (C21).
#   3. Gives the two now-used trailing columns (D, E) explicit widths, and
#      updates the view/selection to the newly active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 3: UserInformation also gets AgeRange ---
$ws.Range("F3").Value = "AgeRange"

# --- Row 20 / new row 21: ProfilePicPath moves from Dashboard-only (C20)
#     down to its own row, now on UserModel (A21) and Dashboard (C21) ---
$ws.Range("C20").Value = $null
$ws.Range("A21").Value = "ProfilePicPath"
$ws.Range("C21").Value = "ProfilePicPath"

# --- Column sizing for the columns now holding data (D, E) ---
$ws.Columns.Item(4).ColumnWidth = 8.71
$ws.Columns.Item(5).ColumnWidth = 17

# --- Scroll/selection: land on the newly added row ---
[void]$ws.Range("D21").Select()
